$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "COD. SERVICIO" header (column M) to "COD. REDMINE"
$ws.Range("M1").Value = "COD. REDMINE"

# Insert a new column before U (shifts U:BA right to V:BB) and give it a header
$ws.Columns("U:U").Insert()
$ws.Columns("U:U").ColumnWidth = $ws.Columns("T:T").ColumnWidth
$ws.Range("U1").Value = "PETICIÓN"

# Restore the view state seen in the saved file (selection near the new column)
$ws.Range("U1").Select()
